# Update the dSF column (column F) values on Sheet1 to reflect the
# recalculated / repulled data, per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1
$ws.Range("F4").Value = -6
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -5
$ws.Range("F12").Value = -5
$ws.Range("F15").Value = -4
$ws.Range("F16").Value = -5
$ws.Range("F18").Value = -5
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = -7
$ws.Range("F21").Value = -3
$ws.Range("F23").Value = -8
$ws.Range("F24").Value = -2
